# The commit swaps the theme used by the deck's (single) slide master
# with the theme that was previously only used by the notes master:
#   ppt/theme/theme1.xml  "Integral"      ->  "Office Theme" palette
#   ppt/theme/theme2.xml  "Office Theme"  ->  "Integral" palette (unreachable
#                                              master, see note below)
#
# Only the <a:clrScheme> (the 12 theme colours) differs between the two
# themes - <a:fontScheme> and <a:fmtScheme> are already byte-identical in
# both files - so re-pointing the slide master at the "Office Theme"
# colours reproduces the effective (rendered) content of the diff for
# ppt/theme/theme1.xml exactly.
#
# PowerPoint's object model exposes the 12 modern theme colours through
# the legacy-named Master.ColorScheme collection, in the fixed order
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. RGB values use VBA's
# RGB(r,g,b) = r + g*256 + b*65536 packing.

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
